$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'https://www.varoom.com/property/beautiful-house-in-mansfield/HA-3214239599'
$ws.Range("F2").Value = 'The property ''Beautiful House in Mansfield'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B3").Value = 'https://www.varoom.com/property/newly-constructed-mansfield-home-with-fenced-yard/BC-8480246'
$ws.Range("F3").Value = 'The property ''Newly Constructed Mansfield Home with Fenced Yard!'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B4").Value = 'https://www.varoom.com/property/newly-constructed-mansfield-home-w-fenced-yard/EP-93205463'
$ws.Range("F4").Value = 'The property ''Newly Constructed Mansfield Home w/Fenced Yard!'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B5").Value = 'https://www.varoom.com/property/stylish-modern-house-in-mansfield/HA-3214239622'
$ws.Range("F5").Value = 'The property ''Stylish Modern House in Mansfield'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B6").Value = 'https://www.varoom.com/property/covered-patio-and-large-yard-mansfield-home/BC-10869407'
$ws.Range("F6").Value = 'The property ''Covered Patio and Large yard Mansfield Home'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B7").Value = 'https://www.varoom.com/property/mansfield-home-w-private-yard-covered-patio/EP-99482582'
$ws.Range("F7").Value = 'The property ''Mansfield Home w/Private Yard & Covered Patio!'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B8").Value = 'https://www.varoom.com/property/beautiful-5bd-3ba-near-bld-waterpark-at-t-stadium/HA-3212396444'
$ws.Range("F8").Value = 'The property ''Beautiful 5BD/3BA Near BLD/Waterpark/AT&T Stadium'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B9").Value = 'https://www.varoom.com/property/private-pool-w-full-amenities-by-bld-joe-pool-lk/BC-6383743'
$ws.Range("F9").Value = 'The property ''PRIVATE POOL w/Full Amenities by BLD & Joe Pool LK'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B10").Value = 'https://www.varoom.com/property/holiday-inn-express-suites-mansfield-an-ihg-hotel/EP-1723202'
$ws.Range("F10").Value = 'The property ''Holiday Inn Express & Suites Mansfield, an IHG Hotel'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B11").Value = 'https://www.varoom.com/property/amazing-views-nature-3-porches-wi-fi-secluded-5mi-to-downtown-dine-shop/HA-3213240954'
$ws.Range("F11").Value = 'The property ''Amazing Views 🕊️Nature 3 Porches ⚡️Wi-fi Secluded .5mi to downtown dine/shop'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B12").Value = 'https://www.varoom.com/property/grand-home-comfort-for-everyone/BC-13224292'
$ws.Range("F12").Value = 'The property ''Grand home - Comfort for everyone'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B13").Value = 'https://www.varoom.com/property/best-western-plus-mansfield-inn-suites/EP-1829967'
$ws.Range("F13").Value = 'The property ''Best Western Plus Mansfield Inn & Suites'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B14").Value = 'https://www.varoom.com/property/newly-constructed-mansfield-home-w-fenced-yard/HA-1219685632'
$ws.Range("F14").Value = 'The property ''Newly Constructed Mansfield Home w/Fenced Yard!'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B15").Value = 'https://www.varoom.com/property/the-estate-home/BC-12251366'
$ws.Range("F15").Value = 'The property ''The Estate Home'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B16").Value = 'https://www.varoom.com/property/fairfield-inn-suites-by-marriott-dallas-mansfield/EP-2516425'
$ws.Range("F16").Value = 'The property ''Fairfield Inn & Suites by Marriott Dallas Mansfield'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B17").Value = 'https://www.varoom.com/property/mansfield-home-w-private-yard-covered-patio/HA-1219878395'
$ws.Range("F17").Value = 'The property ''Mansfield Home w/Private Yard & Covered Patio!'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B18").Value = 'https://www.varoom.com/property/brand-new-smart-home/BC-12224186'
$ws.Range("F18").Value = 'The property ''Brand New Smart Home'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B19").Value = 'https://www.varoom.com/property/comfort-inn-suites-mansfield/EP-807999'
$ws.Range("F19").Value = 'The property ''Comfort Inn & Suites Mansfield'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B20").Value = 'https://www.varoom.com/property/enjoy-resort-style-luxury-accommodations-at-the-estate-home/HA-1217682575'
$ws.Range("F20").Value = 'The property ''ENJOY RESORT STYLE LUXURY ACCOMMODATIONS AT THE ESTATE HOME'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B21").Value = 'https://www.varoom.com/property/holiday-inn-express-hotel-suites-mansfield-an-ihg-hotel/BC-184861'
$ws.Range("F21").Value = 'The property ''Holiday Inn Express Hotel & Suites Mansfield, an IHG Hotel'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B22").Value = 'https://www.varoom.com/property/hampton-inn-suites-mansfield/EP-2246519'
$ws.Range("F22").Value = 'The property ''Hampton Inn & Suites Mansfield'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B23").Value = 'https://www.varoom.com/property/large-comfy-4br-retreat-great-location-fam-group/HA-3213151814'
$ws.Range("F23").Value = 'The property ''Large Comfy 4BR Retreat/Great Location ~Fam/Group'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B24").Value = 'https://www.varoom.com/property/best-western-plus-mansfield-inn-and-suites/BC-254822'
$ws.Range("F24").Value = 'The property ''Best Western Plus Mansfield Inn and Suites'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'

$ws.Range("B25").Value = 'https://www.varoom.com/property/la-quinta-inn-suites-by-wyndham-mansfield-tx/EP-1844878'
$ws.Range("F25").Value = 'The property ''La Quinta Inn & Suites by Wyndham Mansfield TX'' is Available in the specified date range. | **Location:** West Debbieton | **Date Range:** Check-in: 2025-01-25, Check-out: 2025-01-30'
